$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells are treated as text (many contain dot-grouped
# numbers like '25.958.78' or plain decimals that must stay literal text).
$priceCells = @("D2","D3","D5","D10","D11","D12","D13","D16","D17","D18","D20","D23","D25","D28","D31","D34","D35","D36","D40","D42","D43","D45","D46","D47","D51")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated values from the latest cryptos snapshot
$ws.Range("D2").Value = '25.958.78'
$ws.Range("E2").Value = '  +0.63%  '
$ws.Range("D3").Value = '1.587.40'
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '210.26'
$ws.Range("E5").Value = '  +0.51%  '
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -0.37%  '
$ws.Range("E9").Value = '  -0.99%  '
$ws.Range("D10").Value = '17.82'
$ws.Range("E10").Value = '  -1.03%  '
$ws.Range("D11").Value = '0.0808'
$ws.Range("E11").Value = '  +2.39%  '
$ws.Range("D12").Value = '1.809.23'
$ws.Range("E12").Value = '  +0.22%  '
$ws.Range("D13").Value = '1.598.27'
$ws.Range("E13").Value = '  +0.82%  '
$ws.Range("E14").Value = '  -1.30%  '
$ws.Range("E15").Value = '  +0.19%  '
$ws.Range("D16").Value = '25.940.23'
$ws.Range("E16").Value = '  +0.55%  '
$ws.Range("D17").Value = '59.92'
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").Value = '0.0₃0719'
$ws.Range("E18").Value = '  -0.36%  '
$ws.Range("E19").Value = '  -0.21%  '
$ws.Range("D20").Value = '198.43'
$ws.Range("E20").Value = '  +3.61%  '
$ws.Range("E22").Value = '  -1.83%  '
$ws.Range("D23").Value = '5.96'
$ws.Range("E23").Value = '  +0.79%  '
$ws.Range("E24").Value = '  +8.87%  '
$ws.Range("D25").Value = '143.03'
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("E27").Value = '  -8.37%  '
$ws.Range("D28").Value = '15.03'
$ws.Range("E28").Value = '  -0.39%  '
$ws.Range("E29").Value = '  -0.19%  '
$ws.Range("E30").Value = '  +0.25%  '
$ws.Range("D31").Value = '0.0473'
$ws.Range("E31").Value = '  +0.38%  '
$ws.Range("E32").Value = '  +0.22%  '
$ws.Range("E33").Value = '  -2.92%  '
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").Value = '1.47'
$ws.Range("E34").Value = '  -1.69%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '2.37'
$ws.Range("E35").Value = '  +0.59%  '
$ws.Range("D36").Value = '1.122.68'
$ws.Range("E36").Value = '  +2.16%  '
$ws.Range("E37").Value = '  +7.62%  '
$ws.Range("E38").Value = '  -0.14%  '
$ws.Range("E39").Value = '  -1.18%  '
$ws.Range("D40").Value = '0.778'
$ws.Range("E40").Value = '  +0.24%  '
$ws.Range("E41").Value = '  -3.48%  '
$ws.Range("D42").Value = '0.785'
$ws.Range("E42").Value = '  -3.51%  '
$ws.Range("D43").Value = '1.721.29'
$ws.Range("E43").Value = '  +0.14%  '
$ws.Range("E44").Value = '  -2.41%  '
$ws.Range("D45").Value = '91.69'
$ws.Range("E45").Value = '  -2.19%  '
$ws.Range("D46").Value = '1.47'
$ws.Range("E46").Value = '  -1.69%  '
$ws.Range("D47").Value = '53.11'
$ws.Range("E47").Value = '  +0.00%  '
$ws.Range("E48").Value = '  -1.12%  '
$ws.Range("E49").Value = '  -0.29%  '
$ws.Range("E50").Value = '  +0.15%  '
$ws.Range("D51").Value = '0.0₇0914'
$ws.Range("E51").Value = '  -18.38%  '
